$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 updates
$ws.Range("G6").Value = 1.47
$ws.Range("H6").Value = 3.6
$ws.Range("I6").Value = 7.8
$ws.Range("J6").Value = 1.95
$ws.Range("K6").Value = 2.12
$ws.Range("L6").Value = 7.2
$ws.Range("M6").Value = 1.03
$ws.Range("N6").Value = 6.65
$ws.Range("O6").Value = 1.34
$ws.Range("P6").Value = 2.75
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.65
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.52
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.57
$ws.Range("W6").Value = 5.3
$ws.Range("X6").Value = 5.9
$ws.Range("Z6").Value = 9.75
$ws.Range("AA6").Value = 13.5
$ws.Range("AC6").Value = 8
$ws.Range("AD6").Value = 7.4
$ws.Range("AE6").Value = 22
$ws.Range("AF6").Value = 120
$ws.Range("AH6").Value = 16.5
$ws.Range("AI6").Value = 55
$ws.Range("AJ6").Value = 25
$ws.Range("AK6").Value = 250
$ws.Range("AL6").Value = 110
$ws.Range("AM6").Value = 100
$ws.Range("AN6").Value = 3.1
$ws.Range("AO6").Value = 6.7
$ws.Range("AP6").Value = 17.5
$ws.Range("AQ6").Value = 21
$ws.Range("AR6").Value = 55
$ws.Range("AT6").Value = 2.5
$ws.Range("AU6").Value = 8.25
$ws.Range("AV6").Value = 90
$ws.Range("AW6").Value = 8.5
$ws.Range("AX6").Value = 50
$ws.Range("AY6").Value = 50
$ws.Range("AZ6").Value = 400
$ws.Range("BA6").Value = 400

# Row 7 updates
$ws.Range("S7").Value = 1.35
$ws.Range("T7").Value = 2.94
